$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-06 Tuesday" "2025-05-07 Wednesday"

Replace-Text "216×5=1080" "639×9=5751"
Replace-Text "633×5=3165" "283×6=1698"
Replace-Text "890×7=6230" "616×5=3080"
Replace-Text "377×8=3016" "124×4=496"
Replace-Text "121×5=605" "538×8=4304"

Replace-Text "900×3=2700" "766×6=4596"
Replace-Text "441×3=1323" "724×4=2896"
Replace-Text "845×8=6760" "822×5=4110"
Replace-Text "390×4=1560" "151×5=755"
Replace-Text "362×6=2172" "714×3=2142"

Replace-Text "197×4=788" "736×5=3680"
Replace-Text "111×2=222" "761×7=5327"
Replace-Text "588×8=4704" "450×5=2250"
Replace-Text "405×8=3240" "744×7=5208"
Replace-Text "202×7=1414" "687×2=1374"

Replace-Text "310×6=1860" "885×8=7080"
Replace-Text "784×8=6272" "542×8=4336"
Replace-Text "256×8=2048" "191×4=764"
Replace-Text "902×8=7216" "113×7=791"
Replace-Text "693×2=1386" "106×5=530"

Replace-Text "166×8=1328" "850×8=6800"
Replace-Text "379×9=3411" "286×2=572"
Replace-Text "308×3=924" "847×9=7623"
Replace-Text "170×9=1530" "811×5=4055"
Replace-Text "333×2=666" "508×2=1016"
